# Commit: "Ajout 1ere partie Gui" - populate the cost plan (CPC) sheet
# with real values for the cost categories and add the detail/total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 label changes: "Prestation de soins medicaux" -> "Prestation de cours "
$ws.Range("A2").Value = "Prestation de cours "

# Row 3: Couts d'assurance
$ws.Range("B3").Value = "Couts d'assurance"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

# Row 4: Cout de location, espace
$ws.Range("B4").Value = "Cout de location, espace"
$ws.Range("C4").Value = 75
$ws.Range("D4").Value = 75
$ws.Range("E4").Value = 75
$ws.Range("F4").Value = 75
$ws.Range("G4").Value = 75

# Row 5: Cout d'energie, electricite
$ws.Range("B5").Value = "Cout d'energie, électricité"
$ws.Range("C5").Value = 84
$ws.Range("D5").Value = 84
$ws.Range("E5").Value = 84
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 84

# Row 6: Total des couts / Cout d'entretien, with zero values
$ws.Range("A6").Value = "Total des couts"
$ws.Range("B6").Value = "Cout d'entretien"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Row 7: Couts d'assurance (repeated block)
$ws.Range("B7").Value = "Couts d'assurance"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 2

# Row 8: Cout de location, espace
$ws.Range("B8").Value = "Cout de location, espace"
$ws.Range("C8").Value = 75
$ws.Range("D8").Value = 75
$ws.Range("E8").Value = 75
$ws.Range("F8").Value = 75
$ws.Range("G8").Value = 75

# Row 9: Cout d'energie, electricite
$ws.Range("B9").Value = "Cout d'energie, électricité"
$ws.Range("C9").Value = 84
$ws.Range("D9").Value = 84
$ws.Range("E9").Value = 84
$ws.Range("F9").Value = 84
$ws.Range("G9").Value = 84

# Row 10: Cout d'entretien
$ws.Range("B10").Value = "Cout d'entretien"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

# Row 11: Total des couts
$ws.Range("A11").Value = "Total des couts"
$ws.Range("C11").Value = 322
$ws.Range("D11").Value = 322
$ws.Range("E11").Value = 322
$ws.Range("F11").Value = 322
$ws.Range("G11").Value = 322
